# Add data.kind / data.range / data.sheet rows to the "Posting Label" sheet's
# generatedForm block (rows 4-6), right after the existing "Testing from" /
# "Verifying" rows. This mirrors the commit:
#   "Added data.kind, data.range, data.sheet to generatedForm's posting labels"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Posting Label")

# The sheet ships protected (sheet protection, no password) - lift it so the
# new cells can be written, then restore protection at the end.
$ws.Unprotect()

# Column B grows slightly wider to fit the new "data.sheet.57" label.
$ws.Columns.Item(2).ColumnWidth = 12.83

# New label/value pairs (rows 4, 5 and 6).
$ws.Range("B4").Value = "data.kind.57"
$ws.Range("C4").Value = "Manifest for test_dataframe_2_xl"

$ws.Range("B5").Value = "data.range.57"
$ws.Range("C5").Value = "A1:I20"

$ws.Range("B6").Value = "data.sheet.57"
$ws.Range("C6").Value = "Manifest"

# Give the new label cells (column B) the same look as the existing
# "Testing from" / "Verifying" labels above them ...
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4:B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ... and give the new value cells (column C) the same "editable data" look
# used for the posting-label values on the Manifest sheet (light-green fill,
# unlocked cell) rather than the grey "Testing from"/"Verifying" value style.
$wsManifest = $wb.Worksheets.Item("Manifest")
$wsManifest.Range("B2").Copy() | Out-Null
$ws.Range("C4:C6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Restore sheet protection.
$ws.Protect($null, $false, $true, $true, $true)
